# List of commands from Nuk.xlsx - add a second OpenRB response column,
# rename headers, move a couple of values, and remove the "Get dish" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Remove the "Get dish" row (row 13) - table auto-resizes to A1:F24.
$ws.Rows.Item(13).Delete()

# 2) Insert a new column before the old "Comment" column (F), shifting it to G.
$ws.Columns.Item(6).Insert()

# Grow the table to include the freshly inserted column.
$lo.Resize($ws.Range("A1:G24"))

# 3) Header row updates.
$ws.Range("D1").Value = "Expexted response OpenRB1"
$ws.Range("F1").Value = "Expexted result OpenRB2"
$ws.Range("G1").Value = "COMMENT"

# 4) Row 3 (Lift/Lower PetriDish): extend argument/response text with MID/TOP,
#    add the new comment explaining the lift positions.
$ws.Range("C3").Value = "NORMAL/BLOOD/CHOCOLAT/STRG; UP/DOWN/MID/TOP"
$ws.Range("D3").Value = "LIFT UP / LIFT DOWN /LIFT MID /LIFT TOP"
$ws.Range("G3").Value = "UP-> top (right below petri), DOWN-> Bottom, MID-> 1 Petri heigh lower than UP, TOP-> 1 Petri higher than UP"

# 5) Rows 4-5 (Grab/Release Petridish): move GRABBED/RELEASED from
#    "Expexted response OpenRB1" into the new "Expexted result OpenRB2" column.
$ws.Range("D4").Value = "N/A"
$ws.Range("F4").Value = "GRABBED"
$ws.Range("D5").Value = "N/A"
$ws.Range("F5").Value = "RELEASED"

# 6) Move the free-text comments that used to live in column F ("Comment")
#    into the new column G now that the structural insert/move happened.
$ws.Range("G8").Value = "(De)Actiavtes suction and moves gantry"
$ws.Range("G9").Value = "Polar arm needs to move to sample location"
$ws.Range("G10").Value = "Polar arm needs to move to cut area"
$ws.Range("G11").Value = "Polar arm needs to move over petridsih"

Write-Host "done"
